{"js": "// Replace the date line and every \"AxB=\" multiplication prompt in the\n// document, in document order. The mapping below is strictly positional\n// (old[i] -> new[i]); some new values collide with later-in-document old\n// values (e.g. index 10's new text \"18\u00d764=\" equals index 24's old text), so\n// we must not do a global text search-and-replace keyed only on the old\n// string. Instead we walk body.paragraphs (which already yields every\n// paragraph, including ones nested in table cells, in reading order) and\n// replace the Nth non-empty paragraph's text with the Nth replacement.\n\nconst replacements = [\n  \"2026-03-01 Sunday\",\n  \"40\u00d791=\",\n  \"81\u00d712=\",\n  \"69\u00d742=\",\n  \"85\u00d730=\",\n  \"42\u00d765=\",\n  \"20\u00d742=\",\n  \"65\u00d759=\",\n  \"93\u00d785=\",\n  \"86\u00d747=\",\n  \"18\u00d764=\",\n  \"63\u00d795=\",\n  \"70\u00d787=\",\n  \"14\u00d751=\",\n  \"85\u00d763=\",\n  \"94\u00d784=\",\n  \"53\u00d765=\",\n  \"88\u00d780=\",\n  \"50\u00d777=\",\n  \"14\u00d784=\",\n  \"83\u00d744=\",\n  \"36\u00d774=\",\n  \"23\u00d773=\",\n  \"42\u00d753=\",\n  \"48\u00d749=\",\n  \"69\u00d720=\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet replacementIndex = 0;\nfor (const paragraph of paragraphs.items) {\n  if (replacementIndex >= replacements.length) break;\n  if (paragraph.text === \"\") continue;\n  paragraph.getRange().insertText(replacements[replacementIndex], \"Replace\");\n  replacementIndex++;\n}\n\nawait context.sync();\n", "ps1": "# Update the title date and every \"AxB=\" multiplication prompt in the single\n# table. The table has 20 rows x 5 columns, but only rows 1, 5, 10, 15 and 20\n# (1-based COM indexing) carry text; the rows in between are blank spacer\n# rows. We address each content cell by its (row, column) position rather\n# than searching for the old text globally, because some of the new values\n# collide with later old values elsewhere in the document (e.g. \"18x64=\" is\n# simultaneously the new value for one cell and the old value of a\n# different, later cell), which would make a plain find-and-replace-all\n# ambiguous/order-dependent.\n\n$d = $word.ActiveDocument\n\n# Title paragraph with the date.\n$d.Paragraphs.Item(1).Range.Text = \"2026-03-01 Sunday\"\n\n$table = $d.Tables.Item(1)\n\n$contentRows = @(1, 5, 10, 15, 20)\n$newValues = @(\n  @(\"40\u00d791=\", \"81\u00d712=\", \"69\u00d742=\", \"85\u00d730=\", \"42\u00d765=\"),\n  @(\"20\u00d742=\", \"65\u00d759=\", \"93\u00d785=\", \"86\u00d747=\", \"18\u00d764=\"),\n  @(\"63\u00d795=\", \"70\u00d787=\", \"14\u00d751=\", \"85\u00d763=\", \"94\u00d784=\"),\n  @(\"53\u00d765=\", \"88\u00d780=\", \"50\u00d777=\", \"14\u00d784=\", \"83\u00d744=\"),\n  @(\"36\u00d774=\", \"23\u00d773=\", \"42\u00d753=\", \"48\u00d749=\", \"69\u00d720=\")\n)\n\nfor ($i = 0; $i -lt $contentRows.Count; $i++) {\n  $row = $contentRows[$i]\n  $rowValues = $newValues[$i]\n  for ($col = 1; $col -le 5; $col++) {\n    $table.Cell($row, $col).Range.Text = $rowValues[$col - 1]\n  }\n}\n"}
